$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new users. Fill column-by-column (A first for both new rows, then
# B, then C, then D) so new logins/names/surnames/passwords are entered in
# the same sequence as they'd be typed into the sheet.

$ws.Range("A5:A6").NumberFormat = "@"
$ws.Range("D5:D6").NumberFormat = "@"

$ws.Range("A5").Value = "337"
$ws.Range("A6").Value = "Shedl"

$ws.Range("B5").Value = "Елена "
$ws.Range("B6").Value = "МАКСИМ"

$ws.Range("C5").Value = "Шалаева "
$ws.Range("C6").Value = "Вихров"

$ws.Range("D5").Value = "123"
$ws.Range("D6").Value = "12345"
